$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.235119684511858
$ws.Range("C2").Value = 0.5071586931836123
$ws.Range("D2").Value = 0.6973440416867049
$ws.Range("E2").Value = 0.2834699150656448
$ws.Range("G2").Value = 0.002679793540247113
$ws.Range("J2").Value = 0.1467305111030441
$ws.Range("K2").Value = 2.662915578627405
$ws.Range("N2").Value = 5.464709227481421
$ws.Range("B3").Value = 2.175917881542034
$ws.Range("C3").Value = 0.494474703544256
$ws.Range("D3").Value = 0.6892302650019815
$ws.Range("E3").Value = 0.2793667667859339
$ws.Range("G3").Value = 0.002686100997040191
$ws.Range("J3").Value = 0.1439233637187414
$ws.Range("K3").Value = 2.593373454143858
$ws.Range("N3").Value = 5.384164746932356
$ws.Range("B4").Value = 2.140946789185307
$ws.Range("C4").Value = 0.4870057640592051
$ws.Range("D4").Value = 0.684629425232913
$ws.Range("E4").Value = 0.2770098796821969
$ws.Range("G4").Value = 0.002690172732646706
$ws.Range("J4").Value = 0.142289315308183
$ws.Range("K4").Value = 2.552325387472763
$ws.Range("N4").Value = 5.33511702224996
$ws.Range("B5").Value = 2.127040831720166
$ws.Range("C5").Value = 0.4840418579045433
$ws.Range("D5").Value = 0.6828500529240671
$ws.Range("E5").Value = 0.2760901148689712
$ws.Range("G5").Value = 0.002691882207659955
$ws.Range("J5").Value = 0.1416458247208894
$ws.Range("K5").Value = 2.536010961083122
$ws.Range("N5").Value = 5.315230292766643
$ws.Range("B6").Value = 2.124752548317218
$ws.Range("C6").Value = 0.4835545060072945
$ws.Range("D6").Value = 0.6825603489809851
$ws.Range("E6").Value = 0.2759398411507306
$ws.Range("G6").Value = 0.002692169102797261
$ws.Range("J6").Value = 0.141540322857054
$ws.Range("K6").Value = 2.53332684396284
$ws.Range("N6").Value = 5.311934143913533
$ws.Range("B7").Value = 2.140757853745527
$ws.Range("C7").Value = 0.4869654695266377
$ws.Range("D7").Value = 0.6846050416930325
$ws.Range("E7").Value = 0.2769973108932788
$ws.Range("G7").Value = 0.002690195583489399
$ws.Range("J7").Value = 0.1422805464404391
$ws.Range("K7").Value = 2.552103696055752
$ws.Range("N7").Value = 5.334848417711044
$ws.Range("B8").Value = 2.214419605672902
$ws.Range("C8").Value = 0.5027187264784914
$ws.Range("D8").Value = 0.6944671301847052
$ws.Range("E8").Value = 0.2820213270632053
$ws.Range("G8").Value = 0.002681927176011455
$ws.Range("J8").Value = 0.1457439455620673
$ws.Range("K8").Value = 2.638593417502875
$ws.Range("N8").Value = 5.436851981544066
$ws.Range("B9").Value = 2.369901914009176
$ws.Range("C9").Value = 0.5361671787496221
$ws.Range("D9").Value = 0.7168463732124621
$ws.Range("E9").Value = 0.2931710037082595
$ws.Range("G9").Value = 0.002667282731063908
$ws.Range("J9").Value = 0.153252288464941
$ws.Range("K9").Value = 2.821412390055855
$ws.Range("N9").Value = 5.640196880945894
$ws.Range("B10").Value = 2.491001232344388
$ws.Range("C10").Value = 0.5623379720755963
$ws.Range("D10").Value = 0.7351671798355426
$ws.Range("E10").Value = 0.3021669974018195
$ws.Range("G10").Value = 0.00265746848647027
$ws.Range("J10").Value = 0.1592148567572735
$ws.Range("K10").Value = 2.96396003628729
$ws.Range("N10").Value = 5.791752677636907
$ws.Range("B11").Value = 2.547615634486135
$ws.Range("C11").Value = 0.5745989064982382
$ws.Range("D11").Value = 0.7439156432616585
$ws.Range("E11").Value = 0.3064371863876971
$ws.Range("G11").Value = 0.002653206362802084
$ws.Range("J11").Value = 0.1620263775576234
$ws.Range("K11").Value = 3.030635813459924
$ws.Range("N11").Value = 5.861198614032446
$ws.Range("B12").Value = 2.569275811687191
$ws.Range("C12").Value = 0.5792936023886455
$ws.Range("D12").Value = 0.7472884665603203
$ws.Range("E12").Value = 0.3080800010733853
$ws.Range("G12").Value = 0.002651621318107263
$ws.Range("J12").Value = 0.1631054407970254
$ws.Range("K12").Value = 3.056150342386843
$ws.Range("N12").Value = 5.887570418492487
$ws.Range("B13").Value = 2.564601017690052
$ws.Range("C13").Value = 0.5782802045732183
$ws.Range("D13").Value = 0.7465593949843026
$ws.Range("E13").Value = 0.3077250412633603
$ws.Range("G13").Value = 0.00265196140205148
$ws.Range("J13").Value = 0.1628724023220229
$ws.Range("K13").Value = 3.050643464392465
$ws.Range("N13").Value = 5.881887451052876
$ws.Range("B14").Value = 2.54939317804434
$ws.Range("C14").Value = 0.5749841014230697
$ws.Range("D14").Value = 0.7441919232231555
$ws.Range("E14").Value = 0.3065718237130568
$ws.Range("G14").Value = 0.002653075381356026
$ws.Range("J14").Value = 0.1621148632538905
$ws.Range("K14").Value = 3.032729565262628
$ws.Range("N14").Value = 5.863366741591449
$ws.Range("B15").Value = 2.540106851448911
$ws.Range("C15").Value = 0.5729718998612441
$ws.Range("D15").Value = 0.7427496004554257
$ws.Range("E15").Value = 0.3058688091889294
$ws.Range("G15").Value = 0.002653761488002355
$ws.Range("J15").Value = 0.1616527286380744
$ws.Range("K15").Value = 3.021791488417875
$ws.Range("N15").Value = 5.852031983253198
$ws.Range("B16").Value = 2.487332230638231
$ws.Range("C16").Value = 0.5615438995447732
$ws.Range("D16").Value = 0.7346038168891766
$ws.Range("E16").Value = 0.3018915274464931
$ws.Range("G16").Value = 0.002657751084110368
$ws.Range("J16").Value = 0.1590331243700973
$ws.Range("K16").Value = 2.959639672706885
$ws.Range("N16").Value = 5.787224493371809
$ws.Range("B17").Value = 2.455348984803095
$ws.Range("C17").Value = 0.5546247164884051
$ws.Range("D17").Value = 0.7297130217497738
$ws.Range("E17").Value = 0.2994972973264254
$ws.Range("G17").Value = 0.002660250290186249
$ws.Range("J17").Value = 0.1574515730018646
$ws.Range("K17").Value = 2.921982180181374
$ws.Range("N17").Value = 5.74759721862182
$ws.Range("B18").Value = 2.437096544213773
$ws.Range("C18").Value = 0.5506784299757328
$ws.Range("D18").Value = 0.7269389308277141
$ws.Range("E18").Value = 0.2981369282292761
$ws.Range("G18").Value = 0.002661706831440121
$ws.Range("J18").Value = 0.156551226924023
$ws.Range("K18").Value = 2.900494682971214
$ws.Range("N18").Value = 5.724851836578495
$ws.Range("B19").Value = 2.430941159471104
$ws.Range("C19").Value = 0.5493480113882185
$ws.Range("D19").Value = 0.7260063514974604
$ws.Range("E19").Value = 0.2976791975009903
$ws.Range("G19").Value = 0.002662203271101193
$ws.Range("J19").Value = 0.1562479815476792
$ws.Range("K19").Value = 2.893248862548148
$ws.Range("N19").Value = 5.717158686904611
$ws.Range("B20").Value = 2.458738791355017
$ws.Range("C20").Value = 0.5553578096810554
$ws.Range("D20").Value = 0.7302296196185978
$ws.Range("E20").Value = 0.2997504343462509
$ws.Range("G20").Value = 0.00265998227385067
$ws.Range("J20").Value = 0.1576189661526257
$ws.Range("K20").Value = 2.925973053963162
$ws.Range("N20").Value = 5.751810713323721
$ws.Range("B21").Value = 2.553854059586683
$ws.Range("C21").Value = 0.5759508376248164
$ws.Range("D21").Value = 0.7448856762246123
$ws.Range("E21").Value = 0.3069098501164689
$ws.Range("G21").Value = 0.002652747395078586
$ws.Range("J21").Value = 0.1623369787948263
$ws.Range("K21").Value = 3.037984074580777
$ws.Range("N21").Value = 5.86880469574362
$ws.Range("B22").Value = 2.617309673261104
$ws.Range("C22").Value = 0.5897114744071246
$ws.Range("D22").Value = 0.754814003509324
$ws.Range("E22").Value = 0.311739348694843
$ws.Range("G22").Value = 0.002648187521294685
$ws.Range("J22").Value = 0.1655044860562072
$ws.Range("K22").Value = 3.112740678014916
$ws.Range("N22").Value = 5.945700037594975
$ws.Range("B23").Value = 2.583323253550304
$ws.Range("C23").Value = 0.5823393389828482
$ws.Range("D23").Value = 0.7494829301644756
$ws.Range("E23").Value = 0.3091479214341462
$ws.Range("G23").Value = 0.002650605849743597
$ws.Range("J23").Value = 0.1638061907205071
$ws.Range("K23").Value = 3.072698856332067
$ws.Range("N23").Value = 5.904619296642636
$ws.Range("B24").Value = 2.457205839000778
$ws.Range("C24").Value = 0.5550262797292191
$ws.Range("D24").Value = 0.7299959483491136
$ws.Range("E24").Value = 0.2996359409496634
$ws.Range("G24").Value = 0.002660103382562176
$ws.Range("J24").Value = 0.1575432599710638
$ws.Range("K24").Value = 2.924168274609826
$ws.Range("N24").Value = 5.749905677681227
$ws.Range("B25").Value = 2.326644388270267
$ws.Range("C25").Value = 0.5268409510910033
$ws.Range("D25").Value = 0.7104641568035106
$ws.Range("E25").Value = 0.2900144501302151
$ws.Range("G25").Value = 0.002671077618665708
$ws.Range("J25").Value = 0.1511434116912156
$ws.Range("K25").Value = 2.770522556023195
$ws.Range("N25").Value = 5.584817627916834
